$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.248.63'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '2.583.41'
$ws.Range('E3').Value = '  +1.97%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '315.87'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('E7').Value = '  -0.44%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.999'
$ws.Range('E8').Value = '  -0.11%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.541'
$ws.Range('E9').Value = '  +1.03%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '35.61'
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('E12').Value = '  -1.94%  '
$ws.Range('D13').Value = '2.979.66'
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('E14').Value = '  -3.42%  '
$ws.Range('D15').Value = '2.573.42'
$ws.Range('E15').Value = '  +1.53%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '15.16'
$ws.Range('E16').Value = '  -1.53%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.845'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').Value = '43.274.49'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('E19').Value = '  +2.68%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '12.56'
$ws.Range('E20').Value = '  -4.06%  '
$ws.Range('D21').Value = '0.0₃0963'
$ws.Range('E21').Value = '  -0.52%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '69.52'
$ws.Range('E22').Value = '  -1.07%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '254.05'
$ws.Range('E23').Value = '  +0.70%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.98'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  +2.85%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '27.27'
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('E27').Value = '  -0.01%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.45'
$ws.Range('E28').Value = '  +0.61%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '40.26'
$ws.Range('E29').Value = '  +0.62%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '10.34'
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('E31').Value = '  -3.37%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '154.84'
$ws.Range('E32').Value = '  +0.18%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.41'
$ws.Range('E33').Value = '  +2.96%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '2.16'
$ws.Range('E34').Value = '  +1.53%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.0806'
$ws.Range('E35').Value = '  +1.70%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.70'
$ws.Range('E36').Value = '  +3.59%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '18.75'
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('E39').Value = '  +5.49%  '
$ws.Range('E40').Value = '  -0.58%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '22.43'
$ws.Range('E41').Value = '  -5.05%  '
$ws.Range('E42').Value = '  +3.71%  '
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('D46').Value = '2.006.73'
$ws.Range('E46').Value = '  -0.79%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '8.95'
$ws.Range('E47').Value = '  +1.57%  '
$ws.Range('D48').Value = '2.830.58'
$ws.Range('E48').Value = '  +1.85%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '83.10'
$ws.Range('E49').Value = '  -3.30%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '74.99'
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('E51').Value = '  +1.85%  '
